$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "insercao" values between row 2 (KIVO11) and row 5 (FIGS11):
# Row 2 B:D goes from 0,0,0 -> 2,2.5,2.4
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2.5
$ws.Range("D2").Value = 2.4

# Row 5 B:D goes from 2,2.5,2.4 -> 0,0,0
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

# Update the active selection to F5
$ws.Range("F5").Select()
